$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D (Price) and E (Volume) to Text format so that
# numeric-looking strings (e.g. "0.5235", "4.500") are written back
# verbatim as text instead of being auto-converted to numbers by Excel
# (which would silently drop meaningful trailing zeros, change
# "26.224.99" style multi-dot values, etc.).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '26.224.99'
$ws.Cells.Item(2, 5).Value = '  -0.48%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.659.50'
$ws.Cells.Item(3, 5).Value = '  -0.48%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.66%  '

# Row 5
$ws.Cells.Item(5, 5).Value = '  -0.47%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '0.5235'
$ws.Cells.Item(6, 5).Value = '  -2.05%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '1.005'
$ws.Cells.Item(7, 5).Value = '  -0.60%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '0.2638'
$ws.Cells.Item(8, 5).Value = '  -0.99%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.06311'
$ws.Cells.Item(9, 5).Value = '  -1.28%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '20.62'
$ws.Cells.Item(10, 5).Value = '  -1.23%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '0.07793'
$ws.Cells.Item(11, 5).Value = '  -0.58%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '4.500'
$ws.Cells.Item(12, 5).Value = '  -1.26%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.672.36'
$ws.Cells.Item(13, 5).Value = '  +0.78%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '1.887.91'
$ws.Cells.Item(14, 5).Value = '  -0.42%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '0.5636'
$ws.Cells.Item(15, 5).Value = '  +1.64%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '0.0₅8045'
$ws.Cells.Item(16, 5).Value = '  -1.74%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '65.16'
$ws.Cells.Item(17, 5).Value = '  -1.30%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '26.227.08'
$ws.Cells.Item(18, 5).Value = '  -0.59%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '1.003'
$ws.Cells.Item(19, 5).Value = '  -0.86%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '4.712'
$ws.Cells.Item(20, 5).Value = '  +0.64%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '194.14'
$ws.Cells.Item(21, 5).Value = '  -0.69%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '10.21'
$ws.Cells.Item(22, 5).Value = '  -0.53%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '6.024'
$ws.Cells.Item(23, 5).Value = '  -0.34%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '1.006'
$ws.Cells.Item(24, 5).Value = '  -0.66%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '145.47'
$ws.Cells.Item(25, 5).Value = '  -0.60%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '0.1207'
$ws.Cells.Item(26, 5).Value = '  -1.51%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '7.221'
$ws.Cells.Item(27, 5).Value = '  -0.11%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '16.01'
$ws.Cells.Item(28, 5).Value = '  -1.10%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '1.491'
$ws.Cells.Item(29, 5).Value = '  -0.67%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '0.05651'
$ws.Cells.Item(30, 5).Value = '  -3.53%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '1.275'
$ws.Cells.Item(31, 5).Value = '  -0.65%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '3.475'
$ws.Cells.Item(32, 5).Value = '  -3.13%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '3.362'
$ws.Cells.Item(33, 5).Value = '  +2.14%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '1.601'
$ws.Cells.Item(34, 5).Value = '  -0.77%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '2.804'
$ws.Cells.Item(35, 5).Value = '  -1.07%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'HuobiToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(36, 4).Value = '2.405'
$ws.Cells.Item(36, 5).Value = '  -0.71%  '

# Row 37
$ws.Cells.Item(37, 2).Value = 'ARBITRUM'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(37, 4).Value = '0.9425'
$ws.Cells.Item(37, 5).Value = '  -3.00%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '0.5748'
$ws.Cells.Item(38, 5).Value = '  -1.37%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.01601'
$ws.Cells.Item(39, 5).Value = '  -0.52%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '5.998'
$ws.Cells.Item(40, 5).Value = '  +2.36%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '1.054.46'
$ws.Cells.Item(41, 5).Value = '  -2.00%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '2.576'
$ws.Cells.Item(42, 5).Value = '  -0.11%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '0.8456'
$ws.Cells.Item(43, 5).Value = '  -2.18%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.68%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '102.91'
$ws.Cells.Item(45, 5).Value = '  -1.31%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '1.799.20'
$ws.Cells.Item(46, 5).Value = '  -0.35%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '58.12'
$ws.Cells.Item(47, 5).Value = '  +0.03%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  +2.93%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '1.004'
$ws.Cells.Item(49, 5).Value = '  -0.80%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '0.05313'
$ws.Cells.Item(50, 5).Value = '  +2.81%  '

# Row 51
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).Value = '8.035'
$ws.Cells.Item(51, 5).Value = '  -0.34%  '
